$d = $word.ActiveDocument

$replacements = @(
    @{old="44÷5="; new="70÷6="},
    @{old="69÷8="; new="24÷7="},
    @{old="52÷7="; new="27÷2="},
    @{old="96÷7="; new="10÷3="},
    @{old="21÷4="; new="86÷5="},
    @{old="12÷6="; new="81÷3="},
    @{old="49÷8="; new="58÷7="},
    @{old="45÷7="; new="95÷9="},
    @{old="96÷4="; new="66÷8="},
    @{old="41÷4="; new="22÷3="},
    @{old="35÷8="; new="35÷2="},
    @{old="86÷8="; new="31÷9="},
    @{old="19÷5="; new="89÷9="},
    @{old="54÷5="; new="58÷6="},
    @{old="78÷5="; new="32÷7="},
    @{old="87÷7="; new="94÷2="},
    @{old="98÷8="; new="17÷4="},
    @{old="81÷5="; new="80÷6="},
    @{old="23÷8="; new="67÷6="},
    @{old="88÷9="; new="36÷3="},
    @{old="51÷2="; new="70÷4="},
    @{old="34÷9="; new="56÷2="},
    @{old="83÷4="; new="19÷6="},
    @{old="63÷8="; new="32÷6="},
    @{old="33÷5="; new="80÷5="}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $true, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
